$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 300 ("「どこでも自分でありなさい」...") and shift
# everything below it up by one row (matches the dimension shrinking from
# A1:C451 to A1:C450 in the target XML).
$ws.Rows.Item(300).Delete()
